$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Fruta / hortaliza, semanal" -- the weekly refresh reshuffles the Alcachofa
# (artichoke) price records across rows 2-17: each row's per-record fields
# (Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad de
# comercializacion, Precio $/Kg, Kg o Unidades) move to land on a different
# row, while the market/region/category/origin/classification columns
# (A,B,C,E,F,G,O,R) -- and row 15 entirely -- are left untouched.
#
# $rowSource[targetRow] = sourceRow -- i.e. row `targetRow` ends up holding
# the data that row `sourceRow` held before this edit.
$rowSource = @{
    2  = 16
    3  = 12
    4  = 14
    5  = 4
    6  = 17
    7  = 9
    8  = 3
    9  = 2
    10 = 11
    11 = 13
    12 = 6
    13 = 7
    14 = 8
    16 = 10
    17 = 5
}
$cols = @("D", "H", "I", "J", "K", "L", "M", "N", "P", "Q")

# Snapshot every source cell's current value up front (via the Value()
# getter) before any writes happen -- several rows trade data with each
# other, so reading lazily mid-write would clobber a still-needed source.
$orig = @{}
foreach ($r in 2..17) {
    foreach ($c in $cols) {
        $orig["$c$r"] = $ws.Range("$c$r").Value()
    }
}

# Write the permuted values back.
foreach ($targetRow in $rowSource.Keys) {
    $sourceRow = $rowSource[$targetRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $orig["$c$sourceRow"]
    }
}
